# "added device create test"
#
# The "Test Results" sheet tracks, per device, whether each CRUD test
# passed (columns B..E = Create/Read/Update/Delete Test Passed). Rows 2-10
# already record a passing "Create" test (text "True" in column B). This
# adds the same passing "Create" test result for the remaining devices in
# rows 11-24 (column B), leaving the Read/Update/Delete columns untouched.
#
# Column B already stores the Boolean-looking value as literal text ("True"),
# not a native Boolean -- assigning the string "True" directly would be
# auto-coerced to a real Boolean by Excel. Instead, copy the existing text
# value from B2 (which is already the correct literal-text "True") down into
# B11:B24 via PasteSpecial(values), which preserves its text type/format
# exactly like the other "Create Test Passed" cells above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

$source = $ws.Range("B2")
$source.Copy()

for ($r = 11; $r -le 24; $r++) {
    $ws.Range("B$r").PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
